$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating point rounding in A12 (date/time serial value)
$ws.Range("A12").Value = 45866.54189236111

# Append new row 13 with the latest automated reading
$ws.Range("A13").Value = 45866.58354523038
$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 31
$ws.Range("D13").Value = 20.89
$ws.Range("E13").Value = 70.7
$ws.Range("F13").Value = 101.95
$ws.Range("G13").Value = 17.94
$ws.Range("H13").Value = "SE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "14:00:18"

# Match the date/time number format used by the rest of column A
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
